$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-18 12:50:40"
$wsZhCn.Range("G5").Value = "2016-01-18 12:51:23"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-18 12:50:49"
$wsDeDe.Range("G5").Value = "2016-01-18 12:51:41"
